$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1294.8889
$ws.Range("I96").Value = 906.5
$ws.Range("J96").Value = 1605.6
$ws.Range("K96").Value = 2719.5
$ws.Range("L96").Value = 4816.799999999999
$ws.Range("M96").Value = -1346.5
$ws.Range("N96").Value = -7562.799999999999

$ws.Range("H112").Value = 987.28815
$ws.Range("J112").Value = 995.6896400000001
$ws.Range("L112").Value = 2987.06892
$ws.Range("N112").Value = -5203.06892

$ws.Range("H137").Value = 950.9796
$ws.Range("I137").Value = 830.5641000000001
$ws.Range("K137").Value = 2491.6923
$ws.Range("M137").Value = 58.30769999999984

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 5123.6665
$ws.Range("I28").Value = 1685.5
$ws.Range("J28").Value = 12000
$ws.Range("K28").Value = 1685.5
$ws.Range("L28").Value = 12000
$ws.Range("M28").Value = -1493.5
$ws.Range("N28").Value = -12384

$ws.Range("H45").Value = 59813.94
$ws.Range("I45").Value = 100482
$ws.Range("K45").Value = 100482
$ws.Range("M45").Value = -100105

$ws.Range("H61").Value = 5209517.5
$ws.Range("I61").Value = 5748309
$ws.Range("J61").Value = 1200
$ws.Range("K61").Value = 5748309
$ws.Range("L61").Value = 1200
$ws.Range("M61").Value = -5748097
$ws.Range("N61").Value = -1624

$ws.Range("H74").Value = 1368.0312
$ws.Range("I74").Value = 1453.7
$ws.Range("J74").Value = 1225.25
$ws.Range("K74").Value = 1453.7
$ws.Range("L74").Value = 1225.25
$ws.Range("M74").Value = -579.7
$ws.Range("N74").Value = -2973.25

$ws.Range("H77").Value = 1368.0312
$ws.Range("I77").Value = 1453.7
$ws.Range("J77").Value = 1225.25
$ws.Range("K77").Value = 7268.5
$ws.Range("L77").Value = 6126.25
$ws.Range("M77").Value = -2900.5
$ws.Range("N77").Value = -14862.25

$ws.Range("H99").Value = 5123.6665
$ws.Range("I99").Value = 1685.5
$ws.Range("J99").Value = 12000
$ws.Range("K99").Value = 1685.5
$ws.Range("L99").Value = 12000
$ws.Range("M99").Value = 1309.5
$ws.Range("N99").Value = -17990

$ws.Range("H110").Value = 2109.611
$ws.Range("I110").Value = 1356.3846
$ws.Range("J110").Value = 4068
$ws.Range("K110").Value = 1356.3846
$ws.Range("L110").Value = 4068
$ws.Range("M110").Value = 688.6153999999999
$ws.Range("N110").Value = -8158

$ws.Range("H122").Value = 1430.7273
$ws.Range("I122").Value = 675
$ws.Range("J122").Value = 1862.5714
$ws.Range("K122").Value = 2025
$ws.Range("L122").Value = 5587.7142
$ws.Range("M122").Value = 425
$ws.Range("N122").Value = -10487.7142

$ws.Range("H124").Value = 21809.666
$ws.Range("J124").Value = 21809.666
$ws.Range("L124").Value = 21809.666
$ws.Range("N124").Value = -31629.666

$ws.Range("H132").Value = 891.0513
$ws.Range("I132").Value = 760.96875
$ws.Range("J132").Value = 1485.7142
$ws.Range("K132").Value = 2282.90625
$ws.Range("L132").Value = 4457.142599999999
$ws.Range("M132").Value = 247.09375
$ws.Range("N132").Value = -9517.142599999999

$ws.Range("H136").Value = 5209517.5
$ws.Range("I136").Value = 5748309
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 17244927
$ws.Range("L136").Value = 3600
$ws.Range("M136").Value = -17242377
$ws.Range("N136").Value = -8700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3704926.2
$ws.Range("I134").Value = 1143.8636
$ws.Range("J134").Value = 13890328
$ws.Range("K134").Value = 3431.5908
$ws.Range("L134").Value = 41670984
$ws.Range("M134").Value = -896.5907999999999
$ws.Range("N134").Value = -41676054

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1684863
$ws.Range("I31").Value = 2223535.2
$ws.Range("J31").Value = 1512.25
$ws.Range("K31").Value = 2223535.2
$ws.Range("L31").Value = 1512.25
$ws.Range("M31").Value = -2223240.2
$ws.Range("N31").Value = -2102.25

$ws.Range("H34").Value = 1684863
$ws.Range("I34").Value = 2223535.2
$ws.Range("J34").Value = 1512.25
$ws.Range("K34").Value = 2223535.2
$ws.Range("L34").Value = 1512.25
$ws.Range("M34").Value = -2223333.2
$ws.Range("N34").Value = -1916.25

$ws.Range("H58").Value = 30303690
$ws.Range("I58").Value = 35714936
$ws.Range("J58").Value = 720
$ws.Range("K58").Value = 35714936
$ws.Range("L58").Value = 720
$ws.Range("M58").Value = -35714733
$ws.Range("N58").Value = -1126

$ws.Range("H105").Value = 12059.9
$ws.Range("I105").Value = 13377.375
$ws.Range("J105").Value = 6790
$ws.Range("K105").Value = 13377.375
$ws.Range("L105").Value = 6790
$ws.Range("M105").Value = -11630.375
$ws.Range("N105").Value = -10284

$ws.Range("H107").Value = 515.8293
$ws.Range("I107").Value = 361.5357
$ws.Range("K107").Value = 361.5357
$ws.Range("M107").Value = 1558.4643

$ws.Range("H132").Value = 15874610
$ws.Range("I132").Value = 1225.9166
$ws.Range("J132").Value = 37039120
$ws.Range("K132").Value = 3677.7498
$ws.Range("L132").Value = 111117360
$ws.Range("M132").Value = -1147.7498
$ws.Range("N132").Value = -111122420

$ws.Range("H134").Value = 1469.2
$ws.Range("I134").Value = 1374
$ws.Range("J134").Value = 1850
$ws.Range("K134").Value = 4122
$ws.Range("L134").Value = 5550
$ws.Range("M134").Value = -1587
$ws.Range("N134").Value = -10620

$ws.Range("H136").Value = 30303690
$ws.Range("I136").Value = 35714936
$ws.Range("J136").Value = 720
$ws.Range("K136").Value = 107144808
$ws.Range("L136").Value = 2160
$ws.Range("M136").Value = -107142258
$ws.Range("N136").Value = -7260

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 2900
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 4166.6665
$ws.Range("K125").Value = 3000
$ws.Range("L125").Value = 12499.9995
$ws.Range("M125").Value = 1920
$ws.Range("N125").Value = -22339.9995

$ws.Range("H131").Value = 745.39
$ws.Range("J131").Value = 772.6923
$ws.Range("L131").Value = 2318.0769
$ws.Range("N131").Value = -12398.0769

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 12261
$ws.Range("J96").Value = 12261
$ws.Range("L96").Value = 12261
$ws.Range("N96").Value = -17753

$ws.Range("H122").Value = 62513024
$ws.Range("I122").Value = 100020220
$ws.Range("J122").Value = 1033.3334
$ws.Range("K122").Value = 300060660
$ws.Range("L122").Value = 3100.0002
$ws.Range("M122").Value = -300058210
$ws.Range("N122").Value = -8000.0002

$ws.Range("H132").Value = 3280.1333
$ws.Range("I132").Value = 2720.4
$ws.Range("J132").Value = 4399.6
$ws.Range("K132").Value = 8161.200000000001
$ws.Range("L132").Value = 13198.8
$ws.Range("M132").Value = -5631.200000000001
$ws.Range("N132").Value = -18258.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 24397544
$ws.Range("I132").Value = 41668708
$ws.Range("J132").Value = 14724
$ws.Range("K132").Value = 125006124
$ws.Range("L132").Value = 44172
$ws.Range("M132").Value = -125003594
$ws.Range("N132").Value = -49232

$ws.Range("H136").Value = 39410580
$ws.Range("I136").Value = 8405772
$ws.Range("J136").Value = 83334056
$ws.Range("K136").Value = 25217316
$ws.Range("L136").Value = 250002168
$ws.Range("M136").Value = -25214766
$ws.Range("N136").Value = -250007268

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 896.9231
$ws.Range("I81").Value = 805
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 1610
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -549
$ws.Range("N81").Value = -6122

$ws.Range("H84").Value = 896.9231
$ws.Range("I84").Value = 805
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 8050
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -2746
$ws.Range("N84").Value = -30608

$ws.Range("H107").Value = 499.2857
$ws.Range("I107").Value = 495
$ws.Range("K107").Value = 1485
$ws.Range("M107").Value = 435

$ws.Range("H109").Value = 20341.8
$ws.Range("J109").Value = 20341.8
$ws.Range("L109").Value = 20341.8
$ws.Range("N109").Value = -23115.8

$ws.Range("H122").Value = 24444.408
$ws.Range("I122").Value = 28765.389
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 86296.167
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -83846.167
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 33060.605
$ws.Range("I132").Value = 62470.707
$ws.Range("J132").Value = 1812.375
$ws.Range("K132").Value = 187412.121
$ws.Range("L132").Value = 5437.125
$ws.Range("M132").Value = -184882.121
$ws.Range("N132").Value = -10497.125

$ws.Range("H136").Value = 9260920
$ws.Range("I136").Value = 16667623
$ws.Range("J136").Value = 2540.6667
$ws.Range("K136").Value = 50002869
$ws.Range("L136").Value = 7622.000100000001
$ws.Range("M136").Value = -50000319
$ws.Range("N136").Value = -12722.0001
